$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "60.684.32"
$ws.Range("E2").Value = "  -3.89%  "
Set-TextValue $ws.Range("D3") "2.903.52"
$ws.Range("E3").Value = "  -4.25%  "
$ws.Range("E4").Value = "  +0.07%  "
Set-TextValue $ws.Range("D5") "588.91"
$ws.Range("E5").Value = "  -1.11%  "
Set-TextValue $ws.Range("D6") "144.17"
$ws.Range("E6").Value = "  -6.39%  "
$ws.Range("E7").Value = "  +0.05%  "
Set-TextValue $ws.Range("D8") "0.502"
$ws.Range("E8").Value = "  -2.58%  "
Set-TextValue $ws.Range("D9") "2.902.28"
$ws.Range("E9").Value = "  -4.18%  "
Set-TextValue $ws.Range("D10") "6.63"
$ws.Range("E10").Value = "  -4.91%  "
$ws.Range("E11").Value = "  -4.99%  "
Set-TextValue $ws.Range("D12") "0.443"
$ws.Range("E12").Value = "  -4.61%  "
$ws.Range("E13").Value = "  -4.13%  "
Set-TextValue $ws.Range("D14") "33.35"
$ws.Range("E14").Value = "  -6.42%  "
$ws.Range("E15").Value = "  +1.64%  "
Set-TextValue $ws.Range("D16") "3.384.73"
$ws.Range("E16").Value = "  -4.38%  "
Set-TextValue $ws.Range("D17") "60.682.53"
$ws.Range("E17").Value = "  -3.70%  "
$ws.Range("E18").Value = "  -6.04%  "
Set-TextValue $ws.Range("D19") "2.902.60"
$ws.Range("E19").Value = "  -4.30%  "
Set-TextValue $ws.Range("D20") "427.28"
$ws.Range("E20").Value = "  -5.74%  "
Set-TextValue $ws.Range("D21") "13.51"
$ws.Range("E21").Value = "  -5.50%  "
Set-TextValue $ws.Range("D22") "0.683"
$ws.Range("E22").Value = "  -2.07%  "
Set-TextValue $ws.Range("D23") "7.05"
$ws.Range("E23").Value = "  -6.45%  "
Set-TextValue $ws.Range("D24") "81.06"
$ws.Range("E24").Value = "  -2.78%  "
Set-TextValue $ws.Range("D25") "10.76"
$ws.Range("E25").Value = "  -6.16%  "
Set-TextValue $ws.Range("D26") "2.20"
$ws.Range("E26").Value = "  -6.41%  "
Set-TextValue $ws.Range("D27") "11.87"
$ws.Range("E27").Value = "  -4.41%  "
$ws.Range("E28").Value = "  -0.01%  "
Set-TextValue $ws.Range("D29") "2.22"
$ws.Range("E29").Value = "  -3.03%  "
Set-TextValue $ws.Range("D30") "1.00"
$ws.Range("E30").Value = "  +0.14%  "
Set-TextValue $ws.Range("D31") "2.61"
$ws.Range("E31").Value = "  -3.68%  "
Set-TextValue $ws.Range("D32") "7.07"
$ws.Range("E32").Value = "  -6.80%  "
Set-TextValue $ws.Range("D33") "26.37"
$ws.Range("E33").Value = "  -4.66%  "
$ws.Range("E34").Value = "  -4.81%  "
Set-TextValue $ws.Range("D35") "0.0₃0843"
$ws.Range("E35").Value = "  -2.66%  "
$ws.Range("E36").Value = "  -3.17%  "
Set-TextValue $ws.Range("D37") "5.59"
$ws.Range("E37").Value = "  -5.39%  "
Set-TextValue $ws.Range("D38") "3.02"
$ws.Range("E38").Value = "  -5.22%  "
Set-TextValue $ws.Range("D39") "49.27"
$ws.Range("E39").Value = "  -2.49%  "
$ws.Range("E40").Value = "  -6.64%  "
$ws.Range("E41").Value = "  -5.87%  "
Set-TextValue $ws.Range("D42") "8.56"
$ws.Range("E42").Value = "  -6.21%  "
Set-TextValue $ws.Range("D43") "0.291"
$ws.Range("E43").Value = "  -6.11%  "
Set-TextValue $ws.Range("D44") "40.67"
$ws.Range("E44").Value = "  -8.20%  "
$ws.Range("E45").Value = "  -3.69%  "
Set-TextValue $ws.Range("D46") "371.53"
$ws.Range("E46").Value = "  -5.42%  "
Set-TextValue $ws.Range("D47") "2.690.36"
$ws.Range("E47").Value = "  -1.22%  "
Set-TextValue $ws.Range("D48") "132.05"
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("E49").Value = "  +0.04%  "
Set-TextValue $ws.Range("D50") "23.99"
$ws.Range("E50").Value = "  -7.10%  "
$ws.Range("E51").Value = "  -3.32%  "
